# Update the "想去人数" (column F) counts on the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 103
    4  = 1583
    5  = 602
    7  = 3
    8  = 11350
    9  = 17
    10 = 89
    12 = 347
    15 = 12326
    16 = 12979
    18 = 137
    20 = 34
    22 = 44
    23 = 90
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
